# Applies the "Updated cryptos list" data refresh to the Price (D) and
# Volume(1h) (E) columns for rows 2-51, matching the upstream GitHub Actions
# commit that refreshed the scraped coinranking.com figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells look like plain decimals (e.g. "20.50"); format them as
# Text first so Excel keeps the exact original string (incl. trailing zeros)
# instead of silently coercing the assignment to a Double.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.274.17'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '1.662.41'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.74%  '
$ws.Range('D5').Value = '218.53'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').Value = '0.5313'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  +0.71%  '
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('D9').Value = '0.06361'
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('D10').Value = '20.50'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '4.551'
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').Value = '1.671.07'
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('D14').Value = '1.892.34'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').Value = '0.0₅8184'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Value = '65.62'
$ws.Range('E17').Value = '  +0.49%  '
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').Value = '4.656'
$ws.Range('E19').Value = '  +2.41%  '
$ws.Range('D20').Value = '192.45'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').Value = '10.19'
$ws.Range('E21').Value = '  +0.97%  '
$ws.Range('D22').Value = '6.052'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('D23').Value = '1.010'
$ws.Range('D24').Value = '144.77'
$ws.Range('E24').Value = '  +3.02%  '
$ws.Range('E25').Value = '  -1.83%  '
$ws.Range('D26').Value = '7.232'
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').Value = '1.484'
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('D29').Value = '0.05855'
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').Value = '3.583'
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('D32').Value = '3.299'
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('E33').Value = '  +4.30%  '
$ws.Range('D34').Value = '0.9582'
$ws.Range('E34').Value = '  +1.20%  '
$ws.Range('D35').Value = '2.817'
$ws.Range('E35').Value = '  +1.99%  '
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').Value = '0.5805'
$ws.Range('E37').Value = '  +2.93%  '
$ws.Range('D38').Value = '0.01611'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('D39').Value = '5.894'
$ws.Range('E39').Value = '  +0.79%  '
$ws.Range('D40').Value = '0.8531'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('E41').Value = '  +0.68%  '
$ws.Range('D42').Value = '1.047.22'
$ws.Range('E42').Value = '  +3.24%  '
$ws.Range('D43').Value = '104.16'
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('E44').Value = '  +0.33%  '
$ws.Range('D45').Value = '57.25'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('E47').Value = '  +1.03%  '
$ws.Range('D48').Value = '0.4373'
$ws.Range('E48').Value = '  +2.03%  '
$ws.Range('D49').Value = '7.967'
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').Value = '1.441'
$ws.Range('E51').Value = '  -2.05%  '
